$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iteration #3")

# Fill in the two new daily log rows (17 & 18) that were previously blank.
# Column A: date (formatted like the rows above, "d-mmm"), Column B: task
# label (existing shared string), Column C: hours spent.
$ws.Range("A17").NumberFormat = "d-mmm"
$ws.Range("A17").Value = 43220
$ws.Range("B17").Value = "Suite implémentation de la BD dans l'app"
$ws.Range("C17").Value = 4

$ws.Range("A18").NumberFormat = "d-mmm"
$ws.Range("A18").Value = 43221
$ws.Range("B18").Value = "Suite implémentation de la BD dans l'app"
$ws.Range("C18").Value = 4

# Update the active selection on this sheet to reflect where the user was
# working (matches the commit's recorded cursor position).
[void]$ws.Range("C19").Select()
